$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.046.16"
$ws.Range("E2").Value = "  -2.36%  "

# Row 3
$ws.Range("D3").Value = "1.825.84"
$ws.Range("E3").Value = "  -1.15%  "

# Row 4
$ws.Range("E4").Value = "  -0.99%  "

# Row 5
$ws.Range("D5").Value = "'311.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.20%  "

# Row 6
$ws.Range("E6").Value = "  -0.92%  "

# Row 7
$ws.Range("D7").Value = "'0.4233"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.58%  "

# Row 8
$ws.Range("D8").Value = "'0.3670"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.15%  "

# Row 9
$ws.Range("D9").Value = "'0.07222"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.51%  "

# Row 10
$ws.Range("D10").Value = "'0.8438"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.57%  "

# Row 11
$ws.Range("D11").Value = "'20.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.69%  "

# Row 12
$ws.Range("D12").Value = "1.823.31"
$ws.Range("E12").Value = "  -1.24%  "

# Row 13
$ws.Range("D13").Value = "'6.650"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.38%  "

# Row 14
$ws.Range("D14").Value = "'5.291"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.61%  "

# Row 15
$ws.Range("D15").Value = "'0.07038"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.22%  "

# Row 16
$ws.Range("D16").Value = "'89.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.74%  "

# Row 17
$ws.Range("D17").Value = "'1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.90%  "

# Row 18
$ws.Range("D18").Value = "'0.000008742"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.93%  "

# Row 19
$ws.Range("E19").Value = "  -0.92%  "

# Row 20
$ws.Range("D20").Value = "'14.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.56%  "

# Row 21
$ws.Range("D21").Value = "27.092.46"
$ws.Range("E21").Value = "  -2.24%  "

# Row 22
$ws.Range("D22").Value = "'5.132"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.47%  "

# Row 23
$ws.Range("E23").Value = "  -1.94%  "

# Row 24
$ws.Range("D24").Value = "2.049.10"
$ws.Range("E24").Value = "  -1.52%  "

# Row 25
$ws.Range("D25").Value = "'1.980"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "

# Row 26
$ws.Range("D26").Value = "'151.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.50%  "

# Row 27
$ws.Range("D27").Value = "'2.255"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.09%  "

# Row 28
$ws.Range("D28").Value = "'18.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.93%  "

# Row 29
$ws.Range("D29").Value = "'5.249"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.87%  "

# Row 30
$ws.Range("D30").Value = "'116.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.82%  "

# Row 31
$ws.Range("D31").Value = "'0.08700"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.57%  "

# Row 32
$ws.Range("E32").Value = "  -3.84%  "

# Row 33
$ws.Range("D33").Value = "'0.7370"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.17%  "

# Row 34
$ws.Range("E34").Value = "  -0.71%  "

# Row 35
$ws.Range("D35").Value = "'4.421"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.63%  "

# Row 36
$ws.Range("D36").Value = "'1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.05%  "

# Row 37
$ws.Range("D37").Value = "'1.093"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.52%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01944"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.79%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05245"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.74%  "

# Row 40
$ws.Range("D40").Value = "'7.333"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.56%  "

# Row 41
$ws.Range("D41").Value = "'2.873"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.75%  "

# Row 42
$ws.Range("D42").Value = "'0.1686"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.57%  "

# Row 43
$ws.Range("D43").Value = "'0.5065"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.28%  "

# Row 44
$ws.Range("D44").Value = "'8.547"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.69%  "

# Row 45
$ws.Range("D45").Value = "'10.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.77%  "

# Row 46
$ws.Range("D46").Value = "'1.963"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.32%  "

# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.4725"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.70%  "

# Row 48
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'105.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.60%  "

# Row 49
$ws.Range("E49").Value = "  -1.10%  "

# Row 50
$ws.Range("D50").Value = "'0.06327"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.09%  "

# Row 51
$ws.Range("D51").Value = "'1.651"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.16%  "
